# Apply the "added title slide layout + functionization (not final)" edit:
#  1. Drop the last two slides from the deck (sldId 260/261), leaving 4 slides.
#  2. Re-purpose slide 1's text box as the deck's title slide: bump the
#     font size (24pt -> 40pt) and replace the lorem-ipsum body copy with
#     the new placeholder title text.

$p = $ppt.ActivePresentation

# --- 1. Remove the trailing two slides -------------------------------------
# Delete from the end so indices of the slides we keep never shift.
$slideCount = $p.Slides.Count
for ($i = $slideCount; $i -gt 4; $i--) {
    $p.Slides.Item($i).Delete()
}

# --- 2. Update the title slide's text box -----------------------------------
$titleSlide = $p.Slides.Item(1)
$textBox = $titleSlide.Shapes.Item(3)
$textRange = $textBox.TextFrame.TextRange

# Set the font size before the text so the new run picks up sz=4000.
$textRange.Font.Size = 40
$textRange.Text = "Title Slide: Placeholder Name"
